# Update cryptos list values (price and volume columns) with the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.024.17"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "'1.643.32"
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'215.28"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").Value = "'0.5073"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.2579"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "'0.06408"
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").Value = "'0.07733"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "'1.650.95"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "'4.256"
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "'1.868.90"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "'0.5451"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").Value = "'0.0₅7966"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "'63.71"
$ws.Range("E17").Value = "  -3.71%  "
$ws.Range("D18").Value = "'26.002.93"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'205.91"
$ws.Range("E20").Value = "  -5.12%  "
$ws.Range("D21").Value = "'4.359"
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("D22").Value = "'10.00"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "'5.978"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'1.947"
$ws.Range("E25").Value = "  +12.78%  "
$ws.Range("D26").Value = "'142.76"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("D28").Value = "'6.868"
$ws.Range("E28").Value = "  -3.54%  "
$ws.Range("D29").Value = "'15.77"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'0.05033"
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("D31").Value = "'1.236"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").Value = "'3.293"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").Value = "'3.216"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "'1.543"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").Value = "'2.339"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "'0.9113"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").Value = "'2.643"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").Value = "'0.5674"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("D39").Value = "'1.132.26"
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").Value = "'0.01562"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.007"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").Value = "'2.550"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "'5.624"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "'0.8201"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "'99.68"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "'1.781.54"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").Value = "'0.4526"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'55.03"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").Value = "'7.774"
$ws.Range("E51").Value = "  -3.74%  "
